# Update the NATMI TPM results with newly recomputed values.
# (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Bmp6/Bmpr1a -> ECs)
$ws.Range("G2").Value = 17.70643966666667
$ws.Range("H2").Value = 53.119319
$ws.Range("I2").Value = 0.4380235920947999
$ws.Range("J2").Value = 0.4380235920947999
$ws.Range("M2").Value = 6.322177333333333
$ws.Range("N2").Value = 18.966532
$ws.Range("O2").Value = 0.08271011762055308
$ws.Range("P2").Value = 0.08271011762055309
$ws.Range("Q2").Value = 111.9432515146342
$ws.Range("R2").Value = 1007.489263631708
$ws.Range("S2").Value = 0.03622898282273806
$ws.Range("T2").Value = 0.03622898282273807

# Row 3 (ECs -> Bmp6/Bmpr1a -> FAPs)
$ws.Range("G3").Value = 17.70643966666667
$ws.Range("H3").Value = 53.119319
$ws.Range("I3").Value = 0.4380235920947999
$ws.Range("J3").Value = 0.4380235920947999
$ws.Range("O3").Value = 0.5401386314560596
$ws.Range("P3").Value = 0.5401386314560597
$ws.Range("Q3").Value = 731.0456859854684
$ws.Range("R3").Value = 6579.411173869215
$ws.Range("S3").Value = 0.2365934635795525
$ws.Range("T3").Value = 0.2365934635795525

# Row 4 (ECs -> Bmp6/Bmpr1a -> MuSCs)
$ws.Range("G4").Value = 17.70643966666667
$ws.Range("H4").Value = 53.119319
$ws.Range("I4").Value = 0.4380235920947999
$ws.Range("J4").Value = 0.4380235920947999
$ws.Range("M4").Value = 27.73243066666667
$ws.Range("N4").Value = 83.197292
$ws.Range("O4").Value = 0.3628105447549136
$ws.Range("P4").Value = 0.3628105447549136
$ws.Range("Q4").Value = 491.0426104093499
$ws.Range("R4").Value = 4419.383493684149
$ws.Range("S4").Value = 0.1589195780634184
$ws.Range("T4").Value = 0.1589195780634184

# Row 5 (ECs -> Bmp6/Bmpr1a -> Resolving-Mac)
$ws.Range("G5").Value = 17.70643966666667
$ws.Range("H5").Value = 53.119319
$ws.Range("I5").Value = 0.4380235920947999
$ws.Range("J5").Value = 0.4380235920947999
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.096171666666667
$ws.Range("N5").Value = 3.288515
$ws.Range("O5").Value = 0.01434070616847367
$ws.Range("P5").Value = 0.01434070616847367
$ws.Range("Q5").Value = 19.40929748014278
$ws.Range("R5").Value = 174.683677321285
$ws.Range("S5").Value = 0.006281567629090889
$ws.Range("T5").Value = 0.006281567629090889

# Row 6 (FAPs -> Bmp6/Bmpr1a -> ECs)
$ws.Range("G6").Value = 1.617245333333334
$ws.Range("H6").Value = 4.851736000000001
$ws.Range("I6").Value = 0.04000756919748267
$ws.Range("J6").Value = 0.04000756919748267
$ws.Range("M6").Value = 6.322177333333333
$ws.Range("N6").Value = 18.966532
$ws.Range("O6").Value = 0.08271011762055308
$ws.Range("P6").Value = 0.08271011762055309
$ws.Range("Q6").Value = 10.22451178883911
$ws.Range("R6").Value = 92.02060609955201
$ws.Range("S6").Value = 0.003309030754036208
$ws.Range("T6").Value = 0.003309030754036209

# Row 7 (FAPs -> Bmp6/Bmpr1a -> FAPs)
$ws.Range("G7").Value = 1.617245333333334
$ws.Range("H7").Value = 4.851736000000001
$ws.Range("I7").Value = 0.04000756919748267
$ws.Range("J7").Value = 0.04000756919748267
$ws.Range("O7").Value = 0.5401386314560596
$ws.Range("P7").Value = 0.5401386314560597
$ws.Range("Q7").Value = 66.77119999110667
$ws.Range("R7").Value = 600.9407999199601
$ws.Range("S7").Value = 0.02160963367421189
$ws.Range("T7").Value = 0.0216096336742119

# Row 8 (FAPs -> Bmp6/Bmpr1a -> MuSCs)
$ws.Range("G8").Value = 1.617245333333334
$ws.Range("H8").Value = 4.851736000000001
$ws.Range("I8").Value = 0.04000756919748267
$ws.Range("J8").Value = 0.04000756919748267
$ws.Range("M8").Value = 27.73243066666667
$ws.Range("N8").Value = 83.197292
$ws.Range("O8").Value = 0.3628105447549136
$ws.Range("P8").Value = 0.3628105447549136
$ws.Range("Q8").Value = 44.8501440776569
$ws.Range("R8").Value = 403.6512966989121
$ws.Range("S8").Value = 0.01451516797485859
$ws.Range("T8").Value = 0.01451516797485859

# Row 9 (FAPs -> Bmp6/Bmpr1a -> Resolving-Mac)
$ws.Range("G9").Value = 1.617245333333334
$ws.Range("H9").Value = 4.851736000000001
$ws.Range("I9").Value = 0.04000756919748267
$ws.Range("J9").Value = 0.04000756919748267
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.096171666666667
$ws.Range("N9").Value = 3.288515
$ws.Range("O9").Value = 0.01434070616847367
$ws.Range("P9").Value = 0.01434070616847367
$ws.Range("Q9").Value = 1.772778512448889
$ws.Range("R9").Value = 15.95500661204
$ws.Range("S9").Value = 0.0005737367943759767
$ws.Range("T9").Value = 0.0005737367943759767

# Row 10 (MuSCs -> Bmp6/Bmpr1a -> ECs)
$ws.Range("G10").Value = 21.099799
$ws.Range("H10").Value = 63.299397
$ws.Range("I10").Value = 0.5219688387077175
$ws.Range("J10").Value = 0.5219688387077175
$ws.Range("M10").Value = 6.322177333333333
$ws.Range("N10").Value = 18.966532
$ws.Range("O10").Value = 0.08271011762055308
$ws.Range("P10").Value = 0.08271011762055309
$ws.Range("Q10").Value = 133.3966709756893
$ws.Range("R10").Value = 1200.570038781204
$ws.Range("S10").Value = 0.04317210404377881
$ws.Range("T10").Value = 0.04317210404377882

# Row 11 (MuSCs -> Bmp6/Bmpr1a -> FAPs)
$ws.Range("G11").Value = 21.099799
$ws.Range("H11").Value = 63.299397
$ws.Range("I11").Value = 0.5219688387077175
$ws.Range("J11").Value = 0.5219688387077175
$ws.Range("O11").Value = 0.5401386314560596
$ws.Range("P11").Value = 0.5401386314560597
$ws.Range("Q11").Value = 871.147295814005
$ws.Range("R11").Value = 7840.325662326045
$ws.Range("S11").Value = 0.2819355342022952
$ws.Range("T11").Value = 0.2819355342022953

# Row 12 (MuSCs -> Bmp6/Bmpr1a -> MuSCs)
$ws.Range("G12").Value = 21.099799
$ws.Range("H12").Value = 63.299397
$ws.Range("I12").Value = 0.5219688387077175
$ws.Range("J12").Value = 0.5219688387077175
$ws.Range("M12").Value = 27.73243066666667
$ws.Range("N12").Value = 83.197292
$ws.Range("O12").Value = 0.3628105447549136
$ws.Range("P12").Value = 0.3628105447549136
$ws.Range("Q12").Value = 585.1487128481027
$ws.Range("R12").Value = 5266.338415632924
$ws.Range("S12").Value = 0.1893757987166366
$ws.Range("T12").Value = 0.1893757987166366

# Row 13 (MuSCs -> Bmp6/Bmpr1a -> Resolving-Mac)
$ws.Range("G13").Value = 21.099799
$ws.Range("H13").Value = 63.299397
$ws.Range("I13").Value = 0.5219688387077175
$ws.Range("J13").Value = 0.5219688387077175
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.096171666666667
$ws.Range("N13").Value = 3.288515
$ws.Range("O13").Value = 0.01434070616847367
$ws.Range("P13").Value = 0.01434070616847367
$ws.Range("Q13").Value = 23.12900183616167
$ws.Range("R13").Value = 208.161016525455
$ws.Range("S13").Value = 0.0074854017450068
$ws.Range("T13").Value = 0.0074854017450068
